# Updated cryptos list on Thu Oct 17 09:13:54 UTC 2024 with GitHub Actions
# Refresh Price (col D) / Volume(1h) (col E) values, and two coin-row swaps
# (rows 27/28 and rows 49/50) to match the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text parses as a plain number (e.g. "597.02") need
# NumberFormat forced to Text first, otherwise Excel's Value setter coerces
# the string into a float (e.g. 597.01999999999998) and drops precision.
# Switching back to the "Normal" style afterwards keeps the cell's style
# index unchanged (matches the source file, which never sets `s=`) while
# the stored value stays the exact text string.

$ws.Range("D2").Value = "67.398.28"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.624.34"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").Value = "2.623.58"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "3.093.81"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000181"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "67.357.20"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "2.625.42"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "594.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.748.68"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.18%  "
$ws.Range("E36").Value = "  -1.60%  "
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.368"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "0.0₆0290"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.89%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +0.08%  "
